# Cambio formato ds_ric_gas/pwe_anno_prec config start run
#
# - ds_ric_gas_anno_prec / ds_ric_pwr_anno_prec (G2/H2) stop being plain
#   yyyymmdd numbers and become text dates formatted "yyyy-mm-dd HH:mm:ss"
#   (matching the "FORMATO DATA" convention used elsewhere in the sheet).
# - A short explanatory note block is added further down column O.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Foglio1")

# --- New note block under the existing GUIDA text (column O) ---
# Written before the G/H edit below so the shared-string table ends up in
# the same insertion order as the source workbook.
$ws.Range("O12").Value = "NOTE ric_anno_prec_gas/pwr"
$ws.Range("O12").Font.Bold = $true

# --- ds_ric_gas_anno_prec / ds_ric_pwr_anno_prec: numeric -> text date ---
$ws.Range("G1:H2").NumberFormat = "@"
$ws.Range("G2").Value = "2022-02-24 00:00:00"
$ws.Range("H2").Value = "2022-02-24 00:00:00"

$ws.Range("O13").Value = "inserire data snapshot da estrarre in preparazione_dati_inigestion"
$ws.Range("O14").Value = "FORMATO DATA: yyyy-mm-dd HH:mm:ss"

# --- View state (best effort; scroll position / window geometry are host
# session state and may not round-trip through every Excel host) ---
try {
    $excel.ActiveWindow.ScrollColumn = 3
    $excel.ActiveWindow.ScrollRow = 1
} catch {}
$ws.Range("J16").Select()
